# Generate Report for Handback
# Updates the localization-status workbook after a de-de handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    (Overview + both language sheets pick this up since they share text)
#  - zh-cn and de-de detail sheets get their "Latest Target File" /
#    "Latest Handback File" / "Latest Handback DateTime" columns filled in,
#    with a hyperlink added on the newly-populated target-file cell.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6353f3f93a63b16e4b7a1dd7ef7223a8784e4488/e2e/a.md"

# ---- Overview sheet: refresh the status text for both languages ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# ---- zh-cn detail sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlA, "", "", "a.md")

$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-25 08:39:35"

$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-25 08:39:35"

$wsZh.Range("I2:I3").Font.Underline = 2
$wsZh.Range("I2:I3").Font.Color = 15570276

$wsZh.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$wsZh.Columns.Item(10).EntireColumn.AutoFit() | Out-Null

# ---- de-de detail sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlA, "", "", "a.md")

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-25 08:39:42"

$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-25 08:39:42"

$wsDe.Range("I2:I3").Font.Underline = 2
$wsDe.Range("I2:I3").Font.Color = 15570276

$wsDe.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$wsDe.Columns.Item(10).EntireColumn.AutoFit() | Out-Null

# ---- Overview columns E/F autosize to the longer status text ----
$wsOverview.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$wsOverview.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

Write-Host "Handback report generated."
